# Rename Sheet1 -> Login
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Login"

# Header row: UserName / Password
$ws.Range("A1").Value = "UserName"
$ws.Range("B1").Value = "Password"

# Bold header text on a yellow fill, applied as a single combined style
# so the workbook only grows one new cell format (font+fill together).
$headerStyle = $wb.Styles.Add("HeaderStyle")
$headerStyle.Font.Bold = $true
$headerStyle.Interior.Color = 65535
$ws.Range("A1:B1").Style = "HeaderStyle"
$wb.Styles("HeaderStyle").Delete()

# Column widths
$ws.Columns.Item(1).ColumnWidth = 14.333333333333334
$ws.Columns.Item(2).ColumnWidth = 12.833333333333334

# Print setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection as left by the author
$ws.Range("H12").Select()
